# Cotações atualizadas - 2025-12-10
# Append the next daily quote row (row 96) to the sheet, mirroring the
# layout of the existing rows (date serial in column A using the same
# number format as the row above, quote strings in B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 96
$prevRow = $newRow - 1

$ws.Cells.Item($newRow, 1).Value = 46001
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = "22,1541"
$ws.Cells.Item($newRow, 3).Value = "16,0288"
$ws.Cells.Item($newRow, 4).Value = "16,0288"
$ws.Cells.Item($newRow, 5).Value = "16,0288"
